$wb = $excel.ActiveWorkbook

foreach ($ws in $wb.Worksheets) {
    # The workbook's font was re-cached from "Tahoma" to "Calibri" (the
    # theme's actual minor-font) when the file was re-saved by a newer
    # Excel build -- apply the same rename across every cell.
    $ws.Cells.Font.Name = "Calibri"

    # Years 2019-2022 (rows 2-5): growth factor 1 -> 0
    $ws.Range("B2:B5").Value = 0

    # Years 2023-2037 (rows 6-20): replace the shared "=1-0.0113" formula
    # with the equivalent literal value -0.0113
    $ws.Range("B6:B20").Value = -0.0113
}

# Restore per-sheet selection / scroll state to match the saved file:
#   - grw_3u (sheet 1) is now the active/visible tab, scrolled to row 6,
#     with B7:B20 selected (active cell B7)
$ws1 = $wb.Worksheets.Item(1)
$ws1.Activate()
$ws1.Application.ActiveWindow.ScrollRow = 6
$ws1.Range("B7:B20").Select()

#   - grw_mea_dmd (sheet 2): scrolled to row 6, B7:B20 selected
$ws2 = $wb.Worksheets.Item(2)
$ws2.Activate()
$ws2.Application.ActiveWindow.ScrollRow = 6
$ws2.Range("B7:B20").Select()

#   - grw_pea_dmd (sheet 3): no longer the active tab, B7:B20 selected
$ws3 = $wb.Worksheets.Item(3)
$ws3.Activate()
$ws3.Range("B7:B20").Select()

# Leave grw_3u as the active sheet/tab
$ws1.Activate()
